
$wb = $excel.ActiveWorkbook

# --- 1. Rename existing sheet, add the new one right after it ---
$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "Torre Hanoi"
$ws2 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $ws1)
$ws2.Name = "N Reinas"

# --- 2. Header row (reuses existing shared strings "N" / "Tiempo") ---
$ws2.Range("A1").Value = "N"
$ws2.Range("B1").Value = "Tiempo"

# --- 3. Data rows (N = 4..34, Tiempo = measured seconds) ---
$ws2.Cells.Item(2, 1).Value = 4
$ws2.Cells.Item(2, 2).Value = 0.0000012
$ws2.Cells.Item(3, 1).Value = 5
$ws2.Cells.Item(3, 2).Value = 0.0000006
$ws2.Cells.Item(4, 1).Value = 6
$ws2.Cells.Item(4, 2).Value = 0.0000056
$ws2.Cells.Item(5, 1).Value = 7
$ws2.Cells.Item(5, 2).Value = 0.0000013
$ws2.Cells.Item(6, 1).Value = 8
$ws2.Cells.Item(6, 2).Value = 0.0000281
$ws2.Cells.Item(7, 1).Value = 9
$ws2.Cells.Item(7, 2).Value = 0.0000107
$ws2.Cells.Item(8, 1).Value = 10
$ws2.Cells.Item(8, 2).Value = 0.0000346
$ws2.Cells.Item(9, 1).Value = 11
$ws2.Cells.Item(9, 2).Value = 0.0000184
$ws2.Cells.Item(10, 1).Value = 12
$ws2.Cells.Item(10, 2).Value = 0.000126
$ws2.Cells.Item(11, 1).Value = 13
$ws2.Cells.Item(11, 2).Value = 0.0000561
$ws2.Cells.Item(12, 1).Value = 14
$ws2.Cells.Item(12, 2).Value = 0.0011998
$ws2.Cells.Item(13, 1).Value = 15
$ws2.Cells.Item(13, 2).Value = 0.0009722
$ws2.Cells.Item(14, 1).Value = 16
$ws2.Cells.Item(14, 2).Value = 0.0083711
$ws2.Cells.Item(15, 1).Value = 17
$ws2.Cells.Item(15, 2).Value = 0.0049066
$ws2.Cells.Item(16, 1).Value = 18
$ws2.Cells.Item(16, 2).Value = 0.043408
$ws2.Cells.Item(17, 1).Value = 19
$ws2.Cells.Item(17, 2).Value = 0.0028602
$ws2.Cells.Item(18, 1).Value = 20
$ws2.Cells.Item(18, 2).Value = 0.25427
$ws2.Cells.Item(19, 1).Value = 21
$ws2.Cells.Item(19, 2).Value = 0.0112101
$ws2.Cells.Item(20, 1).Value = 22
$ws2.Cells.Item(20, 2).Value = 2.52615
$ws2.Cells.Item(21, 1).Value = 23
$ws2.Cells.Item(21, 2).Value = 0.0390273
$ws2.Cells.Item(22, 1).Value = 24
$ws2.Cells.Item(22, 2).Value = 0.688397
$ws2.Cells.Item(23, 1).Value = 25
$ws2.Cells.Item(23, 2).Value = 0.0854838
$ws2.Cells.Item(24, 1).Value = 26
$ws2.Cells.Item(24, 2).Value = 0.762727
$ws2.Cells.Item(25, 1).Value = 27
$ws2.Cells.Item(25, 2).Value = 0.927587
$ws2.Cells.Item(26, 1).Value = 28
$ws2.Cells.Item(26, 2).Value = 6.54898
$ws2.Cells.Item(27, 1).Value = 29
$ws2.Cells.Item(27, 2).Value = 3.51241
$ws2.Cells.Item(28, 1).Value = 30
$ws2.Cells.Item(28, 2).Value = 141.231
$ws2.Cells.Item(29, 1).Value = 31
$ws2.Cells.Item(29, 2).Value = 34.182
$ws2.Cells.Item(30, 1).Value = 32
$ws2.Cells.Item(30, 2).Value = 242.935
$ws2.Cells.Item(31, 1).Value = 33
$ws2.Cells.Item(31, 2).Value = 410.166
$ws2.Cells.Item(32, 1).Value = 34
$ws2.Cells.Item(32, 2).Value = 7145.59


# --- 4. Formatting: copy the exact styles used on "Torre Hanoi" so fills/borders match ---
# Header (orange fill + full box border) -> style index 3 on Torre Hanoi!A1:B1
$ws1.Range("A1:B1").Copy()
$ws2.Range("A1:B1").PasteSpecial(-4122)

# Data rows 2-31 (light fill + full box border) -> style index 1/2 on Torre Hanoi!A2:B2
$ws1.Range("A2:B2").Copy()
$ws2.Range("A2:B31").PasteSpecial(-4122)

# Last data row (row 32) keeps the fill but the bottom border is removed (continues into row 33)
$ws1.Range("A2:B2").Copy()
$ws2.Range("A32:B32").PasteSpecial(-4122)
$ws2.Range("A32:B32").Borders.Item(9).LineStyle = -4142

# Row 33 - blank trailing row that closes the border box (no fill)
$ws2.Range("A33").Borders.Item(7).LineStyle = 1
$ws2.Range("A33").Borders.Item(8).LineStyle = 1
$ws2.Range("A33").Borders.Item(9).LineStyle = 1
$ws2.Range("B33").Borders.Item(10).LineStyle = 1
$ws2.Range("B33").Borders.Item(8).LineStyle = 1
$ws2.Range("B33").Borders.Item(9).LineStyle = 1

$ws2.Range("G4").Select()

# --- 5. Fix up chart1 (Torre Hanoi) formula references now that the sheet was renamed ---
$co1 = $ws1.ChartObjects().Item(1)
$ser1 = $co1.Chart.SeriesCollection(1)
$ser1.Formula = "=SERIES('Torre Hanoi'!`$B`$1,'Torre Hanoi'!`$A`$2:`$A`$42,'Torre Hanoi'!`$B`$2:`$B`$42,1)"

# --- 6. New chart for N Reinas sheet ---
$co2 = $ws2.ChartObjects().Add(0, 0, 433, 216)
$co2.Chart.ChartType = 4
$co2.Chart.SetSourceData($ws2.Range("B1:B32"))
$ser2 = $co2.Chart.SeriesCollection(1)
$ser2.Formula = "=SERIES('N Reinas'!`$B`$1,'N Reinas'!`$A`$2:`$A`$32,'N Reinas'!`$B`$2:`$B`$32,1)"
$co2.Chart.HasLegend = $false
$co2.Name = "Grafico 2"

# --- 7. Make the new sheet the active tab, like the target workbook ---
$ws2.Select()
